# Remove the explanatory paragraph text about the GameManager/Player
# polling workaround, leaving the (now empty) paragraph with its
# bookmark intact (the old "_GoBack" bookmark stays put).
$d = $word.ActiveDocument

$p = $d.Paragraphs(1)
$r = $p.Range
# Exclude the paragraph mark itself so the paragraph (and its bookmark)
# survive; only the run text is cleared.
$r.MoveEnd(1, -1)
$r.Text = ""
